$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 115.8
$ws.Cells.Item(4, 9).Value = 93
$ws.Cells.Item(4, 11).Value = 93
$ws.Cells.Item(4, 13).Value = 21

$ws.Cells.Item(17, 8).Value = 1184590.2
$ws.Cells.Item(17, 10).Value = 1381851.1
$ws.Cells.Item(17, 12).Value = 4145553.3
$ws.Cells.Item(17, 14).Value = -4145889.3

$ws.Cells.Item(34, 8).Value = 3853.75
$ws.Cells.Item(34, 9).Value = 3853.75
$ws.Cells.Item(34, 11).Value = 3853.75
$ws.Cells.Item(34, 13).Value = -3650.75

$ws.Cells.Item(36, 8).Value = 3853.75
$ws.Cells.Item(36, 9).Value = 3853.75
$ws.Cells.Item(36, 11).Value = 3853.75
$ws.Cells.Item(36, 13).Value = -3138.75

$ws.Cells.Item(80, 8).Value = 829.4286
$ws.Cells.Item(80, 9).Value = 531.4
$ws.Cells.Item(80, 11).Value = 1594.2
$ws.Cells.Item(80, 13).Value = -596.1999999999998

$ws.Cells.Item(83, 8).Value = 829.4286
$ws.Cells.Item(83, 9).Value = 531.4
$ws.Cells.Item(83, 11).Value = 4782.599999999999
$ws.Cells.Item(83, 13).Value = 209.4000000000005

$ws.Cells.Item(98, 8).Value = 2859.4722
$ws.Cells.Item(98, 9).Value = 2443.5
$ws.Cells.Item(98, 10).Value = 4315.375
$ws.Cells.Item(98, 11).Value = 2443.5
$ws.Cells.Item(98, 12).Value = 4315.375
$ws.Cells.Item(98, 13).Value = -945.5
$ws.Cells.Item(98, 14).Value = -7311.375

$ws.Cells.Item(122, 8).Value = 2859.4722
$ws.Cells.Item(122, 9).Value = 2443.5
$ws.Cells.Item(122, 10).Value = 4315.375
$ws.Cells.Item(122, 11).Value = 7330.5
$ws.Cells.Item(122, 12).Value = 12946.125
$ws.Cells.Item(122, 13).Value = -4880.5
$ws.Cells.Item(122, 14).Value = -17846.125

$ws.Cells.Item(137, 8).Value = 49262.76
$ws.Cells.Item(137, 9).Value = 1746.6842
$ws.Cells.Item(137, 10).Value = 500665.5
$ws.Cells.Item(137, 11).Value = 5240.0526
$ws.Cells.Item(137, 12).Value = 1501996.5
$ws.Cells.Item(137, 13).Value = -2690.0526
$ws.Cells.Item(137, 14).Value = -1507096.5

$ws.Cells.Item(138, 8).Value = 5119.6313
$ws.Cells.Item(138, 9).Value = 4892.3335
$ws.Cells.Item(138, 10).Value = 5162.25
$ws.Cells.Item(138, 11).Value = 14677.0005
$ws.Cells.Item(138, 12).Value = 15486.75
$ws.Cells.Item(138, 13).Value = -9537.000499999998
$ws.Cells.Item(138, 14).Value = -25766.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 3022.625
$ws.Cells.Item(2, 9).Value = 2598
$ws.Cells.Item(2, 11).Value = 2598
$ws.Cells.Item(2, 13).Value = -2485

$ws.Cells.Item(50, 8).Value = 1766
$ws.Cells.Item(50, 9).Value = 774
$ws.Cells.Item(50, 11).Value = 774
$ws.Cells.Item(50, 13).Value = -60

$ws.Cells.Item(61, 8).Value = 2179.1538
$ws.Cells.Item(61, 9).Value = 2331.6667
$ws.Cells.Item(61, 10).Value = 1538.6
$ws.Cells.Item(61, 11).Value = 2331.6667
$ws.Cells.Item(61, 12).Value = 1538.6
$ws.Cells.Item(61, 13).Value = -2119.6667
$ws.Cells.Item(61, 14).Value = -1962.6

$ws.Cells.Item(102, 8).Value = 1929.9048
$ws.Cells.Item(102, 9).Value = 1914.0588
$ws.Cells.Item(102, 10).Value = 1997.25
$ws.Cells.Item(102, 11).Value = 1914.0588
$ws.Cells.Item(102, 12).Value = 1997.25
$ws.Cells.Item(102, 13).Value = -292.0588
$ws.Cells.Item(102, 14).Value = -5241.25

$ws.Cells.Item(107, 8).Value = 41576
$ws.Cells.Item(107, 10).Value = 41576
$ws.Cells.Item(107, 12).Value = 41576
$ws.Cells.Item(107, 14).Value = -49256

$ws.Cells.Item(116, 8).Value = 3022.625
$ws.Cells.Item(116, 9).Value = 2598
$ws.Cells.Item(116, 11).Value = 2598
$ws.Cells.Item(116, 13).Value = -304

$ws.Cells.Item(136, 8).Value = 2179.1538
$ws.Cells.Item(136, 9).Value = 2331.6667
$ws.Cells.Item(136, 10).Value = 1538.6
$ws.Cells.Item(136, 11).Value = 6995.000100000001
$ws.Cells.Item(136, 12).Value = 4615.799999999999
$ws.Cells.Item(136, 13).Value = -4445.000100000001
$ws.Cells.Item(136, 14).Value = -9715.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 3022.625
$ws.Cells.Item(3, 9).Value = 2598
$ws.Cells.Item(3, 11).Value = 2598
$ws.Cells.Item(3, 13).Value = -2484

$ws.Cells.Item(20, 8).Value = 6089.6562
$ws.Cells.Item(20, 9).Value = 4632.5654
$ws.Cells.Item(20, 10).Value = 9813.333000000001
$ws.Cells.Item(20, 11).Value = 4632.5654
$ws.Cells.Item(20, 12).Value = 9813.333000000001
$ws.Cells.Item(20, 13).Value = -4385.5654
$ws.Cells.Item(20, 14).Value = -10307.333

$ws.Cells.Item(86, 8).Value = 1617.5483
$ws.Cells.Item(86, 10).Value = 1890.1666
$ws.Cells.Item(86, 12).Value = 1890.1666
$ws.Cells.Item(86, 14).Value = -4136.1666

$ws.Cells.Item(89, 8).Value = 1617.5483
$ws.Cells.Item(89, 10).Value = 1890.1666
$ws.Cells.Item(89, 12).Value = 9450.833000000001
$ws.Cells.Item(89, 14).Value = -20682.833

$ws.Cells.Item(105, 8).Value = 16668302
$ws.Cells.Item(105, 9).Value = 25001852
$ws.Cells.Item(105, 10).Value = 1200
$ws.Cells.Item(105, 11).Value = 25001852
$ws.Cells.Item(105, 12).Value = 1200
$ws.Cells.Item(105, 13).Value = -25000105
$ws.Cells.Item(105, 14).Value = -4694

$ws.Cells.Item(107, 8).Value = 2889.8
$ws.Cells.Item(107, 9).Value = 2000
$ws.Cells.Item(107, 11).Value = 2000
$ws.Cells.Item(107, 13).Value = -80

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 1125.25
$ws.Cells.Item(4, 9).Value = 1125.25
$ws.Cells.Item(4, 11).Value = 1125.25
$ws.Cells.Item(4, 13).Value = -1013.25

$ws.Cells.Item(7, 8).Value = 597.6
$ws.Cells.Item(7, 9).Value = 602.7143
$ws.Cells.Item(7, 11).Value = 602.7143
$ws.Cells.Item(7, 13).Value = -489.7143

$ws.Cells.Item(31, 8).Value = 4599.8
$ws.Cells.Item(31, 10).Value = 4500
$ws.Cells.Item(31, 12).Value = 4500
$ws.Cells.Item(31, 14).Value = -5090

$ws.Cells.Item(34, 8).Value = 4599.8
$ws.Cells.Item(34, 10).Value = 4500
$ws.Cells.Item(34, 12).Value = 4500
$ws.Cells.Item(34, 14).Value = -4904

$ws.Cells.Item(103, 8).Value = 125012570
$ws.Cells.Item(103, 9).Value = 125012570
$ws.Cells.Item(103, 11).Value = 125012570
$ws.Cells.Item(103, 13).Value = -125011398

$ws.Cells.Item(132, 8).Value = 4738
$ws.Cells.Item(132, 9).Value = 4738
$ws.Cells.Item(132, 11).Value = 14214
$ws.Cells.Item(132, 13).Value = -11684

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 95
$ws.Cells.Item(2, 9).Value = 86.111115
$ws.Cells.Item(2, 11).Value = 516.66669
$ws.Cells.Item(2, 13).Value = -403.66669

$ws.Cells.Item(4, 8).Value = 7817858
$ws.Cells.Item(4, 10).Value = 32058452
$ws.Cells.Item(4, 12).Value = 96175356
$ws.Cells.Item(4, 14).Value = -96175580

$ws.Cells.Item(5, 8).Value = 1218.2858
$ws.Cells.Item(5, 9).Value = 684.5
$ws.Cells.Item(5, 11).Value = 2053.5
$ws.Cells.Item(5, 13).Value = -1941.5

$ws.Cells.Item(6, 8).Value = 1950.25
$ws.Cells.Item(6, 9).Value = 1950.25
$ws.Cells.Item(6, 11).Value = 5850.75
$ws.Cells.Item(6, 13).Value = -5737.75

$ws.Cells.Item(7, 8).Value = 91.333336
$ws.Cells.Item(7, 9).Value = 66.666664
$ws.Cells.Item(7, 11).Value = 199.999992
$ws.Cells.Item(7, 13).Value = -87.99999199999999

$ws.Cells.Item(10, 8).Value = 86.57143000000001
$ws.Cells.Item(10, 9).Value = 93.72727
$ws.Cells.Item(10, 11).Value = 281.18181
$ws.Cells.Item(10, 13).Value = -142.18181

$ws.Cells.Item(11, 8).Value = 477.89285
$ws.Cells.Item(11, 10).Value = 439.9
$ws.Cells.Item(11, 12).Value = 1319.7
$ws.Cells.Item(11, 14).Value = -1599.7

$ws.Cells.Item(12, 8).Value = 1022.1667
$ws.Cells.Item(12, 9).Value = 1275.25
$ws.Cells.Item(12, 10).Value = 516
$ws.Cells.Item(12, 11).Value = 3825.75
$ws.Cells.Item(12, 12).Value = 1548
$ws.Cells.Item(12, 13).Value = -3652.75
$ws.Cells.Item(12, 14).Value = -1894

$ws.Cells.Item(13, 8).Value = 604.2
$ws.Cells.Item(13, 9).Value = 595.5
$ws.Cells.Item(13, 10).Value = 639
$ws.Cells.Item(13, 11).Value = 1786.5
$ws.Cells.Item(13, 12).Value = 1917
$ws.Cells.Item(13, 13).Value = -1618.5
$ws.Cells.Item(13, 14).Value = -2253

$ws.Cells.Item(80, 8).Value = 3162.3333
$ws.Cells.Item(80, 9).Value = 2493.5
$ws.Cells.Item(80, 10).Value = 4500
$ws.Cells.Item(80, 11).Value = 7480.5
$ws.Cells.Item(80, 12).Value = 13500
$ws.Cells.Item(80, 13).Value = -6544.5
$ws.Cells.Item(80, 14).Value = -15372

$ws.Cells.Item(83, 8).Value = 3162.3333
$ws.Cells.Item(83, 9).Value = 2493.5
$ws.Cells.Item(83, 10).Value = 4500
$ws.Cells.Item(83, 11).Value = 22441.5
$ws.Cells.Item(83, 12).Value = 40500
$ws.Cells.Item(83, 13).Value = -17761.5
$ws.Cells.Item(83, 14).Value = -49860

$ws.Cells.Item(101, 8).Value = 20000
$ws.Cells.Item(101, 10).Value = 20000
$ws.Cells.Item(101, 12).Value = 60000
$ws.Cells.Item(101, 14).Value = -64868

$ws.Cells.Item(129, 8).Value = 3109.5625
$ws.Cells.Item(129, 9).Value = 1439.75
$ws.Cells.Item(129, 10).Value = 3666.1667
$ws.Cells.Item(129, 11).Value = 4319.25
$ws.Cells.Item(129, 12).Value = 10998.5001
$ws.Cells.Item(129, 13).Value = 680.75
$ws.Cells.Item(129, 14).Value = -20998.5001

$ws.Cells.Item(135, 8).Value = 1218.2858
$ws.Cells.Item(135, 9).Value = 684.5
$ws.Cells.Item(135, 11).Value = 6160.5
$ws.Cells.Item(135, 13).Value = -3625.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(22, 8).Value = 5927.25
$ws.Cells.Item(22, 10).Value = 1429.6666
$ws.Cells.Item(22, 12).Value = 1429.6666
$ws.Cells.Item(22, 14).Value = -2487.6666

$ws.Cells.Item(102, 8).Value = 4174.75
$ws.Cells.Item(102, 9).Value = 4174.75
$ws.Cells.Item(102, 11).Value = 4174.75
$ws.Cells.Item(102, 13).Value = -2552.75

$ws.Cells.Item(129, 8).Value = 0
$ws.Cells.Item(129, 10).Value = 0
$ws.Cells.Item(129, 12).Value = 0
$ws.Cells.Item(129, 14).ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(139, 8).Value = 250357.5
$ws.Cells.Item(139, 10).Value = 250357.5
$ws.Cells.Item(139, 12).Value = 250357.5
$ws.Cells.Item(139, 14).Value = -260637.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 9).Value = 0
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(9, 12).Value = 0
$ws.Cells.Item(9, 13).ClearContents()
$ws.Cells.Item(9, 14).ClearContents()

$ws.Cells.Item(32, 8).Value = 10000
$ws.Cells.Item(32, 9).Value = 10000
$ws.Cells.Item(32, 11).Value = 10000
$ws.Cells.Item(32, 13).Value = -9683

$ws.Cells.Item(132, 8).Value = 14222
$ws.Cells.Item(132, 9).Value = 22119.125
$ws.Cells.Item(132, 10).Value = 3692.5
$ws.Cells.Item(132, 11).Value = 66357.375
$ws.Cells.Item(132, 12).Value = 11077.5
$ws.Cells.Item(132, 13).Value = -63827.375
$ws.Cells.Item(132, 14).Value = -16137.5
